$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column-width bookkeeping: the new layout needs width=16 split across
#    columns F:G (6-7) instead of F:H (6-8), plus a fresh width=16 on column
#    J (10). Deleting+reinserting column H detaches it from the F:H group.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).Delete()
$ws.Columns.Item(8).Insert()
$ws.Columns.Item(10).ColumnWidth = 15.17

# ---------------------------------------------------------------------------
# 2) Header row (row 1): columns A-D stay the same; E-M become the new
#    Typist / Typist QC / Client / Lob / Process / Product Name / State /
#    County / Status layout.
# ---------------------------------------------------------------------------
$header = @("Typist","Typist QC","Client","Lob","Process","Product Name","State","County","Status")
for ($i = 0; $i -lt $header.Length; $i++) {
    $ws.Cells.Item(1, 5 + $i).Value = $header[$i]
}

# ---------------------------------------------------------------------------
# 3) Data rows 2 and 3 - full bulk sample typist / typist_qc rewrite.
# ---------------------------------------------------------------------------
$row2 = @("001CAS","SIPL5316","SIPL5688","SIPL0102","SIPL0103","Cypress Ascendant Services LLC","Title","Search & Typing","Current Owner Search","AL","Shelby","WIP")
$row3 = @("002CAS","SIPL5316","SIPL5688","SIPL0102","SIPL0103","Cypress Ascendant Services LLC","Title","Search & Typing","Commitment Typing","AR","Ashley","WIP")

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $row2[$i]
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value = $row3[$i]
}

# ---------------------------------------------------------------------------
# 4) Formatting for the new OrderID-adjacent typist columns (C2:D3): plain
#    black (non-bold) font, thin border, with the C/D seam's inner edge
#    suppressed on D so the pair reads as one boxed unit.
# ---------------------------------------------------------------------------
$ws.Range("C2:D3").Font.Color = 0
$ws.Range("D2:D3").Borders.Item(7).LineStyle = -4142

# ---------------------------------------------------------------------------
# 5) Restore the last-saved selection.
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
